$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting the existing rows 435-463
# (and all their data) down to 436-464.
$ws.Rows(435).Insert()

# Populate the newly inserted row 435 with the new weekly price entry.
$ws.Range("A435").Value = 5
$ws.Range("B435").Value = "Macroferia Regional de Talca"
$ws.Range("C435").Value = "Maule"
$ws.Range("D435").Value = 45021
$ws.Range("E435").Value = 7
$ws.Range("F435").Value = 100112003
$ws.Range("G435").Value = "Ajo"
$ws.Range("H435").Value = "Chino"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 200
$ws.Range("K435").Value = 18000
$ws.Range("L435").Value = 18000
$ws.Range("M435").Value = 18000
$ws.Range("N435").Value = '$/caja 10 kilos'
$ws.Range("O435").Value = "China"
$ws.Range("P435").Value = 1800
$ws.Range("Q435").Value = 10
$ws.Range("R435").Value = "Hortaliza"
